$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the "Meta description" paragraph that immediately
# follows the title heading ("Play Fever Slot for Free - Review and
# Similar Games" / "Meta description: Play Fever slot for free ...").
# ---------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# Step 2: the italic "image prompt" paragraph at the very end of the
# document gets split in two:
#   - a new bold paragraph carrying the title text is inserted right
#     before it;
#   - its own text is swapped for the new meta-description sentence
#     (formatting/run structure otherwise stays the same: italic).
# ---------------------------------------------------------------------
$imgPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create a feature image for*") {
        $imgPara = $p
        break
    }
}

if ($imgPara -ne $null) {
    # Insert the new paragraph right after the paragraph that currently
    # precedes the image-prompt paragraph - i.e. right before it.
    $prevPara = $imgPara.Previous()
    $prevPara.Range.InsertParagraphAfter()

    $newPara = $prevPara.Next()
    $newPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fever Slot for Free - Review and Similar Games</w:t></w:r></w:p>")

    $old = 'Create a feature image for "Fever" that captures the vibrant and fun atmosphere of disco and the Maya culture in a cartoon style. The image should feature a happy Maya warrior wearing glasses, with disco lights and a dance floor in the background, as well as some of the key symbols from the game such as fruit, bells, and the diamond wild symbol. The cartoon style should be bright and eye-catching, with bold outlines and cheerful colors to match the upbeat mood of the game. The Maya warrior should be depicted with a big smile and a thumbs up, to signal to players that this game is sure to put them in a good mood. The overall impression should be one of fun, excitement, and retro charm, inviting players to step onto the dance floor and spin the reels of Fever.'
    $new = 'Play Fever slot for free and read our review. Find similar games to this modern take on a classic slot game.'
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
